$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(7, 1).Value = 42604.891458333332
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(7, 1).PasteSpecial(-4122)

$ws.Cells.Item(7, 2).Value = "Random"
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 42
$ws.Cells.Item(7, 9).Value = 58
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 33
$ws.Cells.Item(7, 13).Value = 67
